$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 23 new blank rows before row 14 (old row 14 - the thin divider row -
# and everything below it shifts down by 23 rows).
$ws.Range("B14:C36").EntireRow.Insert()

# The newly inserted rows come back with a brand-new (unstyled) format; copy
# the normal "data row" formatting (style used throughout the table, the one
# row 38 - i.e. the row right after our insertion, previously row 15 - still
# carries) onto the 23 new rows so they keep the same look used elsewhere.
$ws.Range("B38:C38").Copy()
$ws.Range("B14:C36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new API entry in the first of the newly inserted rows.
$ws.Cells.Item(14, 2).Value2 = "transaction.delete.master.setBusinessDocument"
$ws.Cells.Item(14, 3).Value2 = "Menghapus Data Dokumen Bisnis"

# Restore the selection to reflect where the edit was made.
$ws.Range("C15").Select()
